$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.486263275146484
$ws.Range("B1").Value = 2.738741397857666
$ws.Range("C1").Value = 6.813967704772949
$ws.Range("D1").Value = 1.748903155326843
$ws.Range("E1").Value = 0.896901547908783
